$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" snapshot refresh (column D) + "Volume(1h)" delta
# refresh (column E). Column D sometimes holds values that look like plain
# numbers (e.g. "97.36", "1.00"); those are prefixed with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# text cells) instead of silently coercing them to numbers and dropping
# significant trailing zeros.

$ws.Range('D2').Value = '42.942.03'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.299.83'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''300.62'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '''97.36'
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').Value = '''35.77'
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').Value = '''0.0788'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '''17.92'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('D15').Value = '2.652.87'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '2.292.58'
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').Value = '42.877.94'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').Value = '''12.79'
$ws.Range('E19').Value = '  -5.48%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '''67.88'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '''240.91'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '''2.15'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('D26').Value = '''2.42'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').Value = '''4.03'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').Value = '''25.40'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').Value = '''165.63'
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('E30').Value = '  -1.34%  '
$ws.Range('D31').Value = '''9.04'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').Value = '''33.05'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('D33').Value = '''4.89'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  -4.27%  '
$ws.Range('D36').Value = '''17.06'
$ws.Range('E36').Value = '  -7.17%  '
$ws.Range('D37').Value = '''2.38'
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('D38').Value = '''0.0687'
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '''2.72'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').Value = '2.018.12'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').Value = '''10.16'
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('D46').Value = '''2.12'
$ws.Range('E46').Value = '  -2.66%  '
$ws.Range('D47').Value = '''17.29'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').Value = '''2.90'
$ws.Range('E49').Value = '  -3.15%  '
$ws.Range('D50').Value = '''53.53'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('D51').Value = '2.519.43'
$ws.Range('E51').Value = '  -0.76%  '
